$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix decimal separators in the text values (comma -> period)
$ws.Range("C3").Value = "– 0.03"
$ws.Range("C12").Value = "– 12.51"
$ws.Range("C13").Value = "– 2.99"
